$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Recorded By" column G had its value "System, dnasr281@gmail.com" swapped
# to "dnasr281@gmail.com, System" across all matching rows.
$col = $ws.Columns.Item(7)
$col.Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System", 1)
